# Restore the original "PRODUCT Implementation Project" training-schedule
# template: fix the mangled "Product...Product" title/heading text and
# replace the placeholder "Product" course branding with the real
# "AI/ML" course branding across all five worksheets.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet 1: Training Schedule Overview
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Schedule Overview")

$ws1.Range("A1").Value = "PRODUCT IMPLEMENTATION PROJECT - TRAINING SCHEDULE"
$ws1.Range("C4").Value = "Enterprise PRODUCT Training Schedule"
$ws1.Range("A7").Value = "TRAINING SCHEDULE SUMMARY"

$ws1.Range("A9").Value = "AI/ML Fundamentals (AI-101)"
$ws1.Range("A10").Value = "AI/ML Platform Overview (AI-102)"
$ws1.Range("A11").Value = "Data Analysis for Business (AI-201)"
$ws1.Range("A12").Value = "Advanced ML Techniques (AI-301)"
$ws1.Range("A13").Value = "MLOps for IT Teams (AI-302)"
$ws1.Range("B13").Value = "ML Engineers, IT"
$ws1.Range("A14").Value = "Model Validation & QA (AI-303)"
$ws1.Range("B14").Value = "ML Engineers, QA"
$ws1.Range("A15").Value = "Executive Overview (AI-401)"
$ws1.Range("A16").Value = "Train-the-Trainer (AI-501)"

$ws1.Range("A18").Value = "TRAINING SCHEDULE STATISTICS"

# -----------------------------------------------------------------------
# Sheet 2: Detailed Training Schedule
# -----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Detailed Training Schedule")

$ws2.Range("A1").Value = "DETAILED TRAINING SCHEDULE"

$ws2.Range("A4").Value = "AI-101"
$ws2.Range("B4").Value = "AI/ML Fundamentals"
$ws2.Range("A5").Value = "AI-102"
$ws2.Range("B5").Value = "AI/ML Platform Overview"
$ws2.Range("A6").Value = "AI-201"
$ws2.Range("A7").Value = "AI-201"
$ws2.Range("A8").Value = "AI-201"
$ws2.Range("A9").Value = "AI-301"
$ws2.Range("B9").Value = "Advanced ML Techniques"
$ws2.Range("A10").Value = "AI-301"
$ws2.Range("B10").Value = "Advanced ML Techniques"
$ws2.Range("A11").Value = "AI-302"
$ws2.Range("B11").Value = "MLOps for IT Teams"
$ws2.Range("C11").Value = "ML Engineers, IT"
$ws2.Range("A12").Value = "AI-302"
$ws2.Range("B12").Value = "MLOps for IT Teams"
$ws2.Range("C12").Value = "ML Engineers, IT"
$ws2.Range("A13").Value = "AI-303"
$ws2.Range("C13").Value = "ML Engineers, QA"
$ws2.Range("A14").Value = "AI-303"
$ws2.Range("C14").Value = "ML Engineers, QA"
$ws2.Range("A15").Value = "AI-401"
$ws2.Range("A16").Value = "AI-501"
$ws2.Range("A17").Value = "AI-501"
$ws2.Range("A18").Value = "AI-501"
$ws2.Range("A19").Value = "AI-501"
$ws2.Range("A20").Value = "AI-501"

# -----------------------------------------------------------------------
# Sheet 3: Instructor Schedule
# -----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Instructor Schedule")

$ws3.Range("B7").Value = "Advanced ML Techniques"
$ws3.Range("B8").Value = "Advanced ML Techniques"
$ws3.Range("B9").Value = "MLOps for IT Teams"
$ws3.Range("B10").Value = "MLOps for IT Teams"

# -----------------------------------------------------------------------
# Sheet 4: Facility Schedule
# -----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Facility Schedule")

$ws4.Range("B7").Value = "Advanced ML Techniques"
$ws4.Range("B8").Value = "Advanced ML Techniques"
$ws4.Range("B9").Value = "MLOps for IT Teams"
$ws4.Range("B10").Value = "MLOps for IT Teams"

# -----------------------------------------------------------------------
# Sheet 5: Participant Tracking
# -----------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Participant Tracking")

$ws5.Range("F4").Value = "AI-101"
$ws5.Range("F5").Value = "AI-102"
$ws5.Range("F6").Value = "AI-301"
$ws5.Range("F7").Value = "AI-302"
$ws5.Range("E8").Value = "ML Engineer"
$ws5.Range("F8").Value = "AI-101"
$ws5.Range("E9").Value = "ML Engineer"
$ws5.Range("F9").Value = "AI-102"
$ws5.Range("E10").Value = "ML Engineer"
$ws5.Range("F10").Value = "AI-302"
$ws5.Range("E11").Value = "ML Engineer"
$ws5.Range("F11").Value = "AI-303"
$ws5.Range("F12").Value = "AI-101"
$ws5.Range("F13").Value = "AI-102"
$ws5.Range("F14").Value = "AI-401"
$ws5.Range("F15").Value = "AI-101"
$ws5.Range("F16").Value = "AI-102"
$ws5.Range("F17").Value = "AI-501"
$ws5.Range("F18").Value = "AI-101"
$ws5.Range("F19").Value = "AI-102"
$ws5.Range("F20").Value = "AI-301"
$ws5.Range("F21").Value = "AI-303"
$ws5.Range("F22").Value = "AI-501"
